$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename header L1, fill in previously-missing numeric cells (G/I columns)
#    on the existing rows, and append the new row 7.
# ---------------------------------------------------------------------------
$ws.Range("L1").Value = "MasterSheet RowNo."

$ws.Range("G2").Value = 2
$ws.Range("I2").Value = 0

$ws.Range("G3").Value = 0
$ws.Range("I3").Value = 0

$ws.Range("G4").Value = 0
$ws.Range("I4").Value = 0

$ws.Range("G5").Value = 0
$ws.Range("I5").Value = 0

$ws.Range("G6").Value = 2
$ws.Range("I6").Value = 0

$ws.Range("A7").Value = 71
$ws.Range("B7").Value = 241
$ws.Range("C7").Value = "LATIN AMER. & CARIB    "
$ws.Range("D7").Value = 43921
$ws.Range("E7").Value = "Saint Kitts and Nevis"
$ws.Range("F7").Value = 2
$ws.Range("G7").Value = -2
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = "Imported cases only"
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 5368

# ---------------------------------------------------------------------------
# 2) Column widths: A:O => 27 characters wide. Columns M:O have no values but
#    still need to end up inside the used range / dimension (A1:O7); the
#    format-paste pass below takes care of stamping them into existence.

#    (ColumnWidth uses Excel's "characters" unit which has a +5/6 offset vs.
#    the raw OOXML width attribute, so subtract it to land exactly on 27.)
# ---------------------------------------------------------------------------
$ws.Range("A1:O1").EntireColumn.ColumnWidth = 26.1666666666667

# ---------------------------------------------------------------------------
# 4) Styling: build the two new cell formats on scratch cells (kept on the
#    same sheet so Copy/PasteSpecial resolves), then stamp them onto the
#    real ranges one contiguous block at a time. Using single format-only
#    paste operations (rather than chaining multiple format property writes
#    across a multi-cell range) avoids leaving unused intermediate styles
#    behind in the saved style table.
# ---------------------------------------------------------------------------
$scratchGeneral = $ws.Range("Z1")
$scratchGeneral.HorizontalAlignment = -4108
$scratchGeneral.VerticalAlignment = -4108

$scratchDate = $ws.Range("Z2")
$scratchDate.NumberFormat = "yyyy-mm-dd;"
$scratchDate.HorizontalAlignment = -4108
$scratchDate.VerticalAlignment = -4108

# Centered (general) style -> everything except column D.
$scratchGeneral.Copy()
$ws.Range("A1:C1").PasteSpecial(-4122)
$scratchGeneral.Copy()
$ws.Range("E1:O1").PasteSpecial(-4122)
$scratchGeneral.Copy()
$ws.Range("A2:C7").PasteSpecial(-4122)
$scratchGeneral.Copy()
$ws.Range("E2:O7").PasteSpecial(-4122)

# Centered + custom date format -> column D.
$scratchDate.Copy()
$ws.Range("D1:D7").PasteSpecial(-4122)

$ws.Range("Z1:Z2").Clear()
